$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-strings table gained two new entries ("line7", "line8") inserted
# right after "line6" / before "extr1". Because the original file referenced
# the "extr*" labels by their (now-shifted) shared-string index rather than
# re-binding to the same text, every row from 8 downward is relabeled two
# slots further down the name sequence (line7, line8, extr1..extr8), and two
# brand new rows (16, 17) are appended for extr7/extr8.

# Row 8: name -> line7, C 5->14, D 12->11, E stays false
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: name -> line8, C 5->16, D stays 9, E 0->1 (true)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10: name -> extr1, C 10->5, D 11->12, E stays true
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11: name -> extr2, C 7->5, D 8->9, E 0->1 (true)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12: name -> extr3, C 9->10, D stays 11, E 0->1 (true)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

# Row 13: name -> extr4, C stays 7, D 11->8, E 1->0 (false)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14: name -> extr5, C 5->9, D 7->11, E stays false
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15: name -> extr6, C 8->7, D 5->11, E stays true
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# --- Append two new data rows (16, 17), matching the existing row style ---

# Seed new rows by copying formatting from the last existing row, then set values
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A15").Copy($ws.Range("A17"))

# Row 16: index 14, name "extr7", C=5, D=7, E=false
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17: index 15, name "extr8", C=8, D=5, E=true
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
